# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Mon Aug  7 15:38:15 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Worksheet, $Address, $NewValue)
    $rng = $Worksheet.Range($Address)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = $origStyle
}

Set-TextCell $ws 'D2' '29.086.71'
Set-TextCell $ws 'E2' '  +0.02%  '
Set-TextCell $ws 'D3' '1.834.74'
Set-TextCell $ws 'E3' '  +0.10%  '
Set-TextCell $ws 'D4' '1.010'
Set-TextCell $ws 'E4' '  +0.97%  '
Set-TextCell $ws 'D5' '243.13'
Set-TextCell $ws 'E5' '  -0.23%  '
Set-TextCell $ws 'D6' '0.6179'
Set-TextCell $ws 'E6' '  -1.65%  '
Set-TextCell $ws 'E7' '  +0.65%  '
Set-TextCell $ws 'D8' '0.07406'
Set-TextCell $ws 'E8' '  -0.93%  '
Set-TextCell $ws 'D9' '0.2907'
Set-TextCell $ws 'E9' '  -0.55%  '
Set-TextCell $ws 'D10' '22.88'
Set-TextCell $ws 'E10' '  -0.95%  '
Set-TextCell $ws 'D11' '0.07707'
Set-TextCell $ws 'E11' '  -0.19%  '
Set-TextCell $ws 'D12' '1.843.85'
Set-TextCell $ws 'E12' '  -0.34%  '
Set-TextCell $ws 'D13' '4.964'
Set-TextCell $ws 'E13' '  -0.20%  '
Set-TextCell $ws 'D14' '0.6674'
Set-TextCell $ws 'E14' '  -0.16%  '
Set-TextCell $ws 'D15' '82.25'
Set-TextCell $ws 'E15' '  -0.45%  '
Set-TextCell $ws 'D16' '0.000009063'
Set-TextCell $ws 'E16' '  -3.25%  '
Set-TextCell $ws 'D17' '5.866'
Set-TextCell $ws 'E17' '  -2.80%  '
Set-TextCell $ws 'D18' '29.131.56'
Set-TextCell $ws 'E18' '  +0.03%  '
Set-TextCell $ws 'D19' '2.082.98'
Set-TextCell $ws 'E19' '  -0.04%  '
Set-TextCell $ws 'D20' '235.08'
Set-TextCell $ws 'E20' '  +5.48%  '
Set-TextCell $ws 'E21' '  -0.42%  '
Set-TextCell $ws 'E22' '  +0.51%  '
Set-TextCell $ws 'D23' '7.140'
Set-TextCell $ws 'E23' '  +0.15%  '
Set-TextCell $ws 'D24' '1.012'
Set-TextCell $ws 'E24' '  +1.01%  '
Set-TextCell $ws 'D25' '159.07'
Set-TextCell $ws 'E25' '  -0.67%  '
Set-TextCell $ws 'D26' '0.1416'
Set-TextCell $ws 'E26' '  +1.30%  '
Set-TextCell $ws 'D27' '8.475'
Set-TextCell $ws 'E27' '  -0.27%  '
Set-TextCell $ws 'D28' '17.73'
Set-TextCell $ws 'E28' '  -1.03%  '
Set-TextCell $ws 'D29' '1.498'
Set-TextCell $ws 'E29' '  -0.45%  '
Set-TextCell $ws 'B30' 'Hedera'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D30' '0.05558'
Set-TextCell $ws 'E30' '  -2.23%  '
Set-TextCell $ws 'B31' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D31' '4.103'
Set-TextCell $ws 'E31' '  +0.80%  '
Set-TextCell $ws 'D32' '4.112'
Set-TextCell $ws 'E32' '  -0.81%  '
Set-TextCell $ws 'E33' '  +0.74%  '
Set-TextCell $ws 'D34' '1.841'
Set-TextCell $ws 'E34' '  -0.63%  '
Set-TextCell $ws 'D35' '0.7394'
Set-TextCell $ws 'E35' '  -1.25%  '
Set-TextCell $ws 'D36' '1.135'
Set-TextCell $ws 'E36' '  -0.21%  '
Set-TextCell $ws 'D37' '2.652'
Set-TextCell $ws 'E37' '  +1.52%  '
Set-TextCell $ws 'D38' '2.821'
Set-TextCell $ws 'E38' '  +2.54%  '
Set-TextCell $ws 'D39' '0.01772'
Set-TextCell $ws 'E39' '  -0.59%  '
Set-TextCell $ws 'D40' '1.201.93'
Set-TextCell $ws 'E40' '  -2.23%  '
Set-TextCell $ws 'D41' '6.401'
Set-TextCell $ws 'E41' '  -2.33%  '
Set-TextCell $ws 'D42' '0.9083'
Set-TextCell $ws 'E42' '  +1.49%  '
Set-TextCell $ws 'E43' '  +0.56%  '
Set-TextCell $ws 'D44' '101.04'
Set-TextCell $ws 'E44' '  -0.98%  '
Set-TextCell $ws 'D45' '1.976.33'
Set-TextCell $ws 'E45' '  -0.45%  '
Set-TextCell $ws 'D46' '64.74'
Set-TextCell $ws 'E46' '  -1.25%  '
Set-TextCell $ws 'B47' 'Mantle'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell $ws 'D47' '0.5133'
Set-TextCell $ws 'E47' '  +0.95%  '
Set-TextCell $ws 'B48' 'BabyDogeCoin'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 'D48' '0.00000000121'
Set-TextCell $ws 'E48' '  -3.60%  '
Set-TextCell $ws 'D49' '0.4013'
Set-TextCell $ws 'E49' '  -1.07%  '
Set-TextCell $ws 'D50' '9.044'
Set-TextCell $ws 'E50' '  +0.36%  '
Set-TextCell $ws 'E51' '  +0.40%  '
